$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing rows down to make room for new content -------------
# (matches: row4->6, row5->11, row7->15, row8->16, row9->17, row10->18, row12->27)
$ws.Rows("3:4").Insert()
$ws.Rows("7:10").Insert()
$ws.Rows("13:14").Insert()
$ws.Rows("19:25").Insert()

# --- Row 2: weapon motor - update price note ----------------------------
$ws.Range("D2").Value = "6249+1500+600"

# --- Row 6 (was row 4): weapon ESC - Redbrick 200A ----------------------
# (values already shifted down automatically by the row insert above)

# --- Rows 7-10: new ESC sourcing alternatives ----------------------------
$ws.Range("D7").Value = 5499
$ws.Range("E7").Value = "store trustable or not?"
$ws.Range("E8").Value = "trustable or not, also pre-order"
$ws.Range("E9").Value = "says that no post tax - but aroung 80 dollars"

# --- Row 11 (was row 5): Spintend 100A escs - price changed -------------
$ws.Range("D11").Value = 1200

# --- Row 16 (was row 8): Tattu G-Tech battery ----------------------------
$ws.Range("D16").Value = 9999
$ws.Range("E16").Value = "should it be separate batteries for drive and weapon - if yes then calculation for each"

# --- Row 17 (was row 9): Bonka battery -----------------------------------
$ws.Range("E17").Value = "has high C + low C options - weapon + drive we could do"

# --- Row 18 (was row 10): new battery option 3 ---------------------------
$ws.Range("B18").Value = "150c 1050mAh"

# --- Rows 19-24: new battery sourcing options ----------------------------
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "95c 1550mAh"
$ws.Range("E19").Value = "much cheaper than others - but 95 C instead of 150C, still high C - though C rating always overshot"

$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Genx 6S for Drive"
$ws.Range("E20").Value = "lower C - useful for drive"

$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "ChinaHobbyLine"
$ws.Range("E21").Value = "internet says that best FPV batteries for the price - though would have to import"

$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Ovonic"
$ws.Range("E22").Value = "ships to india but only less than 100Wh"

$ws.Range("A24").Value = 8
$ws.Range("B24").Value = "Indian Robo Store"
$ws.Range("E24").Value = "donno how reliable is the store"

# --- Hyperlinks (also paints the "Hyperlink" cell style) -----------------
$ws.Hyperlinks.Add($ws.Range("C16"), "https://www.technobotix.in/products/tattu-g-tech-3500mah-150c-22-2v-6s-lipo-battery/1781252000009101158") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), "https://robu.in/product-category/batteries/batteries-batteries/lithium-polymer-battery-packs/bonka-li-po-battery/6-cell-22-2v-25-2v-bonka-li-po-battery-lithium-polymer-battery-packs/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C18"), "https://www.quadkart.in/tattu-r-line-version-5-0-1050mah-6s-22-2v-150c-lipo-battery-xt60/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C19"), "https://robokits.co.in/batteries-chargers/drone-batteries/genx-power-premium-lipo-battery/genxpower-22.2v-lipo-batteries/genx-22.2v-6s-1550mah-95c-190c-premium-lipo-lithium-polymer-battery") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C20"), "https://robokits.co.in/genxpower-22.2v-lipo-batteries") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C21"), "https://chinahobbyline.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C22"), "https://www.ovonicshop.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C23"), "https://www.moglix.com/brands/ovonic") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C24"), "https://indianrobostore.com/category/drone-batteries/lipo-battery") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://probots.co.in/red-brick-200a-bldc-esc-electronic-speed-controller-2-7s-bec-5v-5a.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.technobotix.in/products/hobbyking-red-brick-200a-esc-v2/1781252000000063739") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.roboverse.in/product-page/red-brick-200a") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "https://rcdrone.top/products/red-brick-speed-controller?variant=43812668113120") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.technobotix.in/products/bbox-pluto-h600-550kv-bldc-motor/1781252000001751636") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.desertcart.in/search?query=red+brick+esc") | Out-Null

# --- Row height for the "Battery - 6S Lipo" header row -------------------
$ws.Rows(15).RowHeight = 30

# --- Column widths (best effort) ------------------------------------------
$ws.Columns("B").ColumnWidth = 51.21875
$ws.Columns("C").ColumnWidth = 35.44140625

# --- Selection -------------------------------------------------------------
$ws.Range("C13").Select() | Out-Null
